$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for rows 2-12 from serial date 45204 to 45207
$ws.Range("C2:C12").Value = 45207
